# Creation of fsp admin panel
# Renames headers, adds a new "delivered_quantity" column (K), restyles the
# new/changed header + data cells, trims the sheet down to the 2 real data
# rows, and resizes the newly-relevant columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename existing headers -------------------------------------------------
$ws.Range("D1").Value = "admin_level_2"
$ws.Range("E1").Value = "collector_name"
$ws.Range("G1").Value = "fsp_name"
$ws.Range("I1").Value = "entitlement_quantity"
$ws.Range("J1").Value = "entitlement_quantity_usd"

# --- 2. Build the formatting for the new K2/K3 (empty data) cells from the ------
#        about-to-be-removed J4 filler cell, which already carries the right
#        fill/border combo (fillId=2 borderId=4) used for the tail columns.
$ws.Range("J4").Copy()
$ws.Range("K2:K3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. New header cell K1: "delivered_quantity", styled like G1/J1 -------------
$ws.Range("G1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K1").Value = "delivered_quantity"

# --- 4. J1 gets its own distinct border (left/right grey, top/bottom red) -------
$rngJ1 = $ws.Range("J1")
$rngJ1.Borders.Item(7).LineStyle = 1
$rngJ1.Borders.Item(7).Color = 10066329
$rngJ1.Borders.Item(10).LineStyle = 1
$rngJ1.Borders.Item(10).Color = 10066329
$rngJ1.Borders.Item(8).LineStyle = 1
$rngJ1.Borders.Item(8).Color = 11184810
$rngJ1.Borders.Item(9).LineStyle = 1
$rngJ1.Borders.Item(9).Color = 11184810

# --- 5. Drop the filler rows 4-10, only the header + 2 real rows remain ---------
$ws.Range("A4:K10").EntireRow.Delete()

# --- 6. Resize the columns that now hold the quantity/fsp data -----------------
$ws.Columns.Item(9).ColumnWidth = 29.5859
$ws.Columns.Item(10).ColumnWidth = 36.6875
$ws.Columns.Item(11).ColumnWidth = 36.6875

Write-Host "ok"
